$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B32").Value = "Notes:"
$ws.Range("C32").Value = "CPWG feedline impedance of 49.95ohm generated with 91mil trace width and 10mil ground spacing on each side."
$ws.Range("C33").Value = "One row of vias placed 20mils from the keepout, 118mils apart from each other, and 61mils away from the CPWG."
